# Refatorando o consolidador para modelo ETL
# Atualiza os dados da planilha de absenteísmo com o novo conjunto gerado pelo ETL.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Linha 2
$ws.Cells.Item(2, 1).Value = 15433
$ws.Cells.Item(2, 2).Value = "Renan Lima"
$ws.Cells.Item(2, 3).Value = "Atendimento ao Cliente"
$ws.Cells.Item(2, 5).Value = 5
$ws.Cells.Item(2, 6).Value = 45095
$ws.Cells.Item(2, 7).Value = 7445.46

# Linha 3
$ws.Cells.Item(3, 1).Value = 33779
$ws.Cells.Item(3, 2).Value = "Pedro Almeida"
$ws.Cells.Item(3, 3).Value = "Marketing"
$ws.Cells.Item(3, 4).Value = "Problemas pessoais"
$ws.Cells.Item(3, 5).Value = 8
$ws.Cells.Item(3, 6).Value = 45087
$ws.Cells.Item(3, 7).Value = 3636.51

# Linha 4
$ws.Cells.Item(4, 1).Value = 10647
$ws.Cells.Item(4, 2).Value = "Dr. Levi Rodrigues"
$ws.Cells.Item(4, 6).Value = 45080
$ws.Cells.Item(4, 7).Value = 10316.87

# Linha 5
$ws.Cells.Item(5, 1).Value = 83456
$ws.Cells.Item(5, 2).Value = "Sr. Luiz Miguel Rezende"
$ws.Cells.Item(5, 3).Value = "TI"
$ws.Cells.Item(5, 4).Value = "Consulta médica"
$ws.Cells.Item(5, 5).Value = 5
$ws.Cells.Item(5, 6).Value = 45091
$ws.Cells.Item(5, 7).Value = 7807.04

# Linha 6
$ws.Cells.Item(6, 1).Value = 50626
$ws.Cells.Item(6, 2).Value = "Dr. João Guilherme da Costa"
$ws.Cells.Item(6, 3).Value = "Vendas"
$ws.Cells.Item(6, 4).Value = "Outros"
$ws.Cells.Item(6, 5).Value = 6
$ws.Cells.Item(6, 6).Value = 45093
$ws.Cells.Item(6, 7).Value = 11729.1

# Linha 7
$ws.Cells.Item(7, 1).Value = 93842
$ws.Cells.Item(7, 2).Value = "João Gabriel da Rocha"
$ws.Cells.Item(7, 3).Value = "TI"
$ws.Cells.Item(7, 4).Value = "Consulta médica"
$ws.Cells.Item(7, 5).Value = 1
$ws.Cells.Item(7, 6).Value = 45093
$ws.Cells.Item(7, 7).Value = 3284.1

# Linha 8
$ws.Cells.Item(8, 1).Value = 7541
$ws.Cells.Item(8, 2).Value = "Lívia Nunes"
$ws.Cells.Item(8, 3).Value = "Financeiro"
$ws.Cells.Item(8, 5).Value = 4
$ws.Cells.Item(8, 6).Value = 45106
$ws.Cells.Item(8, 7).Value = 5632.56

# Linha 9
$ws.Cells.Item(9, 1).Value = 37682
$ws.Cells.Item(9, 2).Value = "Enrico Novaes"
$ws.Cells.Item(9, 3).Value = "Jurídico"
$ws.Cells.Item(9, 5).Value = 3
$ws.Cells.Item(9, 6).Value = 45086
$ws.Cells.Item(9, 7).Value = 6535.47

# Linha 10
$ws.Cells.Item(10, 1).Value = 52226
$ws.Cells.Item(10, 2).Value = "Sr. Eduardo Cavalcanti"
$ws.Cells.Item(10, 3).Value = "Financeiro"
$ws.Cells.Item(10, 5).Value = 3
$ws.Cells.Item(10, 6).Value = 45090
$ws.Cells.Item(10, 7).Value = 5110.23

# Linha 11
$ws.Cells.Item(11, 1).Value = 87
$ws.Cells.Item(11, 2).Value = "Sophie Carvalho"
$ws.Cells.Item(11, 3).Value = "Financeiro"
$ws.Cells.Item(11, 4).Value = "Problemas pessoais"
$ws.Cells.Item(11, 5).Value = 5
$ws.Cells.Item(11, 6).Value = 45099
$ws.Cells.Item(11, 7).Value = 10756.76
